$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 (total_registros) and B2 (errores_fecha) hold plain integer-looking
# text ("8.535" / "6"); force text formatting first so Excel keeps them
# as strings instead of re-interpreting them as numbers.
$ws.Range("A2:B2").NumberFormat = "@"

$ws.Range("A2").Value = "8.535"
$ws.Range("B2").Value = "6"
$ws.Range("C2").Value = "0,07"
$ws.Range("G2").Value = "10,33"
